$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds an alternating March/April daily case-count time series.
# This update ("22 April 1st update") adds the missing entries for
# 21/04/2020 and 22/04/2020 in their correct chronological slots,
# pushing every row below down by two, and also corrects the
# previously-reported value for 20/04/2020.

# Correct the value for 20/04/2020 (row 42).
$ws.Range("B42").Value = 1239

# Insert "21/04/2020" -> 1537 right after "21/03/2020" (row 43).
$ws.Rows.Item(44).Insert()
$ws.Range("A44").Value = "21/04/2020"
$ws.Range("B44").Value = 1537

# Insert "22/04/2020" -> 34 right after "22/03/2020" (now row 45).
$ws.Rows.Item(46).Insert()
$ws.Range("A46").Value = "22/04/2020"
$ws.Range("B46").Value = 34
